$d = $word.ActiveDocument
$apos = [char]0x2019
$ellipsis = [char]0x2026

# ---------------------------------------------------------------------------
# Paragraph "Authentication/authorisation is implemented..." :
#   - unit test count 72 -> 74, drop its yellow highlight
#   - Jacoco test coverage 93% -> 97%, drop the yellow highlight on
#     "Jacoco" / " Test Report"
#   - sonarqube coverage "…" -> "at 92%"
#   - move the "config files are excluded" sentence to the end, covering
#     both reports
# ---------------------------------------------------------------------------

# Replace the highlighted "72" with "74".
$rngNum = $d.Content
$foundNum = $rngNum.Find.Execute("72", $false, $false, $false, $false, $false, $true, 1, $false, "74", 2)

# Clear the (now stray) yellow highlighting left on that run.
$rngClear = $d.Content
$rngClear.Find.Execute("74", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngClear.HighlightColorIndex = 0

# "...93%. The config files are excluded from this report, because I didn't
# make the code for it. The " -> "...97%. The " (stop right before the
# "sonarqube" spell-checked word so its proofErr markers stay intact).
$oldCoverage = "93%. The config files are excluded from this report, because I didn" + $apos + "t make the code for it. The "
$newCoverage = "97%. The "
$rngCoverage = $d.Content
$rngCoverage.Find.Execute($oldCoverage, $false, $false, $false, $false, $false, $false, 1, $false, $newCoverage, 2)

# " coverage is … " -> " coverage is at 92%. The config files are excluded
# from both of these reports, because I didn't make the code for it. "
$oldSonar = " coverage is " + $ellipsis + " "
$newSonar = " coverage is at 92%. The config files are excluded from both of these reports, because I didn" + $apos + "t make the code for it. "
$rngSonar = $d.Content
$rngSonar.Find.Execute($oldSonar, $false, $false, $false, $false, $false, $false, 1, $false, $newSonar, 2)

# ---------------------------------------------------------------------------
# Paragraph "Sonarqube is added to the CI/CD pipeline..." : add a trailing
# space run at the end of the paragraph.
# ---------------------------------------------------------------------------
$pSonarPipeline = $d.Paragraphs.Item(7)
$rngPipeline = $pSonarPipeline.Range
$rngPipeline.MoveEnd(1, -1) | Out-Null
$rngPipeline.InsertAfter(" ")

# ---------------------------------------------------------------------------
# Paragraph "Docker is also added to the CI/CD pipeline..." : explain the
# websockets issue instead of the generic "for some reason".
# ---------------------------------------------------------------------------
$oldDocker = " for some reason."
$newDocker = ", I think because of the websockets because this is the part where it gets stuck."
$rngDocker = $d.Content
$rngDocker.Find.Execute($oldDocker, $false, $false, $false, $false, $false, $false, 1, $false, $newDocker, 2)
